# Apply the "bare essentials for creating work" change:
#  - Defs sheet: replace sample rows 2-3 with real data, add two new rows (4-5)
#  - Point Defs sheet: replace sample rows 2-3 with real data, drop old row 4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Defs" (sheet1.xml) -> dimension grows from A1:I3 to A1:I5
# ---------------------------------------------------------------------------
$defs = $wb.Worksheets.Item("Defs")

$defsRows = @(
    @("'lgqy9rb3-0bmn", "'2023-04-21T14:34:06.217", "'lgqy9rbd", "'FALSE", "'0m7w", "test one",  "1️⃣", "Initial desc",            "SECOND"),
    @("'lgqy9rbd-avpc", "'2023-04-21T14:34:06.217", "'lgqy9rbe", "'FALSE", "'ay7l", "twooo",     "2️⃣", "now with a description",  "WEEK"),
    @("'lgqy9rbe-2ban", "'2023-04-21T14:34:06.218", "'lgqy9rbe", "'FALSE", "'05a8", "afree",     "3️⃣", "Set a description",       "SECOND"),
    @("'lgqy9rbe-0keb", "'2023-04-21T14:34:06.218", "'lgqy9rbe", "'FALSE", "'e0bq", "FOUR",      "4️⃣", "having fun",              "SECOND")
)

$r = 2
foreach ($row in $defsRows) {
    $c = 1
    foreach ($val in $row) {
        $defs.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "Point Defs" (sheet2.xml) -> dimension shrinks from A1:L4 to A1:L3
# ---------------------------------------------------------------------------
$points = $wb.Worksheets.Item("Point Defs")

$pointRows = @(
    @("'lgqy9rbe-0bcq", "'2023-04-21T14:34:06.218", "'lgqy9rbe", $false, "'e0bq", "'0pc6", "set alternatively", "☝️", "Set a description", "BOOL", "COUNT", "TEXT"),
    @("'lgqy9rbe-3tnn", "'2023-04-21T14:34:06.218", "'lgqy9rbe", $false, "'e0bq", "'0tb7", "test point",        "🆕", "Set a description", "TEXT", "COUNT", "TEXT")
)

$r = 2
foreach ($row in $pointRows) {
    $c = 1
    foreach ($val in $row) {
        $points.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# Drop the old 4th data row entirely (was lep65g3sq.fipe / momm / First time? ...)
$points.Range("A4:L4").EntireRow.Delete()
